$d = $word.ActiveDocument

# The HTML->docx import split "let url = "";" across several runs
# (one run per "word", thanks to the spell-checker's proofErr markers):
#   "      " | [spellStart]"let"[spellEnd] | " " | [spellStart]"url"[spellEnd] | " = "";"
# The edit collapses the " " + "url" + " = " portion back into a single
# run and removes the now-redundant proofErr wrapper around "url",
# leaving:
#   "      " | [spellStart]"let"[spellEnd] | " url = "";"
# Locate that paragraph by its text (rather than a hard-coded index) so
# the script is resilient to unrelated structural changes elsewhere in
# the document.
$targetRange = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "*let url = `"`";*") {
        $targetRange = $para.Range
    }
}

if ($targetRange -ne $null) {
    # Only the " url = " substring is searched/replaced - deliberately
    # excluding the literal `"`" characters from the Find/Replace text
    # keeps Word's smart-quote AutoCorrect from mangling them, and
    # scoping the Find to this single paragraph's Range keeps the other
    # (unrelated) "... url = "https://..." assignments further down in
    # the same script block untouched.
    $targetRange.Find.ClearFormatting()
    $targetRange.Find.Execute(" url = ", $false, $false, $false, $false, $false, $true, 1, $false, " url = ", 2)
}
